# Auto-generated Excel COM-interop script to apply the Aegis_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

# ALC!row 112
$ALC.Range("H112").Value = 1336.32
$ALC.Range("J112").Value = 1489.8889
$ALC.Range("L112").Value = 4469.6667
$ALC.Range("N112").Value = -6685.6667

# ALC!row 129
$ALC.Range("H129").Value = 2255.0532
$ALC.Range("I129").Value = 4975.727
$ALC.Range("J129").Value = 1125.717
$ALC.Range("K129").Value = 14927.181
$ALC.Range("L129").Value = 3377.151
$ALC.Range("M129").Value = -9927.181
$ALC.Range("N129").Value = -13377.151

# ALC!row 137
$ALC.Range("H137").Value = 1035.4912
$ALC.Range("I137").Value = 1020.3461
$ALC.Range("J137").Value = 1193
$ALC.Range("K137").Value = 3061.0383
$ALC.Range("L137").Value = 3579
$ALC.Range("M137").Value = -511.0383000000002
$ALC.Range("N137").Value = -8679

# ARM!row 21
$ARM.Range("H21").Value = 7254.2856
$ARM.Range("I21").Value = 2196
$ARM.Range("J21").Value = 19900
$ARM.Range("K21").Value = 2196
$ARM.Range("L21").Value = 19900
$ARM.Range("M21").Value = -1822
$ARM.Range("N21").Value = -20648

# ARM!row 32
$ARM.Range("H32").Value = 5300.91
$ARM.Range("I32").Value = 4190.6445
$ARM.Range("J32").Value = 15293.3
$ARM.Range("K32").Value = 4190.6445
$ARM.Range("L32").Value = 15293.3
$ARM.Range("M32").Value = -3903.6445
$ARM.Range("N32").Value = -15867.3

# ARM!row 74
$ARM.Range("H74").Value = 503.77142
$ARM.Range("I74").Value = 455.0645
$ARM.Range("J74").Value = 881.25
$ARM.Range("K74").Value = 455.0645
$ARM.Range("L74").Value = 881.25
$ARM.Range("M74").Value = 418.9355
$ARM.Range("N74").Value = -2629.25

# ARM!row 77
$ARM.Range("H77").Value = 503.77142
$ARM.Range("I77").Value = 455.0645
$ARM.Range("J77").Value = 881.25
$ARM.Range("K77").Value = 2275.3225
$ARM.Range("L77").Value = 4406.25
$ARM.Range("M77").Value = 2092.6775
$ARM.Range("N77").Value = -13142.25

# BSM!row 94
$BSM.Range("H94").Value = 448.43332
$BSM.Range("I94").Value = 377.33334
$BSM.Range("J94").Value = 732.8333
$BSM.Range("K94").Value = 377.33334
$BSM.Range("L94").Value = 732.8333
$BSM.Range("M94").Value = 73.66665999999998
$BSM.Range("N94").Value = -1634.8333

# BSM!row 134
$BSM.Range("H134").Value = 2225.7593
$BSM.Range("I134").Value = 2029.3125
$BSM.Range("J134").Value = 3797.3333
$BSM.Range("K134").Value = 6087.9375
$BSM.Range("L134").Value = 11391.9999
$BSM.Range("M134").Value = -3552.9375
$BSM.Range("N134").Value = -16461.9999

# CRP!row 31
$CRP.Range("H31").Value = 27634.678
$CRP.Range("I31").Value = 1636.7
$CRP.Range("K31").Value = 1636.7
$CRP.Range("M31").Value = -1341.7

# CRP!row 34
$CRP.Range("H34").Value = 27634.678
$CRP.Range("I34").Value = 1636.7
$CRP.Range("K34").Value = 1636.7
$CRP.Range("M34").Value = -1434.7

# CRP!row 58
$CRP.Range("H58").Value = 1349.1915
$CRP.Range("I58").Value = 1202.6316
$CRP.Range("J58").Value = 1968
$CRP.Range("K58").Value = 1202.6316
$CRP.Range("L58").Value = 1968
$CRP.Range("M58").Value = -999.6315999999999
$CRP.Range("N58").Value = -2374

# CRP!row 132
$CRP.Range("H132").Value = 3058.64
$CRP.Range("I132").Value = 2903.9429
$CRP.Range("K132").Value = 8711.8287
$CRP.Range("M132").Value = -6181.8287

# CRP!row 134
$CRP.Range("H134").Value = 1154.1666
$CRP.Range("I134").Value = 1049.8422
$CRP.Range("J134").Value = 1334.3636
$CRP.Range("K134").Value = 3149.5266
$CRP.Range("L134").Value = 4003.0908
$CRP.Range("M134").Value = -614.5266000000001
$CRP.Range("N134").Value = -9073.0908

# CRP!row 136
$CRP.Range("H136").Value = 1349.1915
$CRP.Range("I136").Value = 1202.6316
$CRP.Range("J136").Value = 1968
$CRP.Range("K136").Value = 3607.8948
$CRP.Range("L136").Value = 5904
$CRP.Range("M136").Value = -1057.8948
$CRP.Range("N136").Value = -11004

# CUL!row 123
$CUL.Range("H123").Value = 3575.1428
$CUL.Range("J123").Value = 3999.2
$CUL.Range("L123").Value = 11997.6
$CUL.Range("N123").Value = -16897.6

# CUL!row 129
$CUL.Range("H129").Value = 286596.78
$CUL.Range("I129").Value = 11975.2
$CUL.Range("J129").Value = 392220.47
$CUL.Range("K129").Value = 35925.60000000001
$CUL.Range("L129").Value = 1176661.41
$CUL.Range("M129").Value = -30925.60000000001
$CUL.Range("N129").Value = -1186661.41

# CUL!row 131
$CUL.Range("H131").Value = 7661.1147
$CUL.Range("J131").Value = 7692.058
$CUL.Range("L131").Value = 23076.174
$CUL.Range("N131").Value = -33156.174

# CUL!row 133
$CUL.Range("H133").Value = 4982.517
$CUL.Range("I133").Value = 3410
$CUL.Range("J133").Value = 5392.7393
$CUL.Range("K133").Value = 10230
$CUL.Range("L133").Value = 16178.2179
$CUL.Range("M133").Value = -5170
$CUL.Range("N133").Value = -26298.2179

# CUL!row 136
$CUL.Range("H136").Value = 1340
$CUL.Range("I136").Value = 1175
$CUL.Range("K136").Value = 3525
$CUL.Range("M136").Value = 1575

# CUL!row 137
$CUL.Range("H137").Value = 3892009
$CUL.Range("I137").Value = 116537.78
$CUL.Range("J137").Value = 5890788
$CUL.Range("K137").Value = 349613.34
$CUL.Range("L137").Value = 17672364
$CUL.Range("M137").Value = -344513.34
$CUL.Range("N137").Value = -17682564

# CUL!row 138
$CUL.Range("H138").Value = 11304.454
$CUL.Range("I138").Value = 14418.625
$CUL.Range("K138").Value = 43255.875
$CUL.Range("M138").Value = -38115.875

# CUL!row 139
$CUL.Range("H139").Value = 2351.1035
$CUL.Range("I139").Value = 1736.4286
$CUL.Range("K139").Value = 5209.2858
$CUL.Range("M139").Value = -69.28579999999965

# CUL!row 140
$CUL.Range("H140").Value = 5623.731
$CUL.Range("I140").Value = 9140.385
$CUL.Range("J140").Value = 2107.077
$CUL.Range("K140").Value = 27421.155
$CUL.Range("L140").Value = 6321.231000000001
$CUL.Range("M140").Value = -22241.155
$CUL.Range("N140").Value = -16681.231

# CUL!row 141
$CUL.Range("H141").Value = 7264.684
$CUL.Range("I141").Value = 8152.0713
$CUL.Range("K141").Value = 24456.2139
$CUL.Range("M141").Value = -19276.2139

# LTW!row 132
$LTW.Range("H132").Value = 2560.6597
$LTW.Range("I132").Value = 2506.25
$LTW.Range("K132").Value = 7518.75
$LTW.Range("M132").Value = -4988.75

# LTW!row 136
$LTW.Range("H136").Value = 1238.4131
$LTW.Range("I136").Value = 984.075
$LTW.Range("J136").Value = 2934
$LTW.Range("K136").Value = 2952.225
$LTW.Range("L136").Value = 8802
$LTW.Range("M136").Value = -402.2250000000004
$LTW.Range("N136").Value = -13902

# WVR!row 28
$WVR.Range("H28").Value = 13257.143
$WVR.Range("I28").Value = 0
$WVR.Range("J28").Value = 13257.143
$WVR.Range("K28").Value = 0
$WVR.Range("L28").Value = 13257.143
$WVR.Range("M28").ClearContents()
$WVR.Range("N28").Value = -13953.143

# WVR!row 132
$WVR.Range("H132").Value = 1720.6323
$WVR.Range("I132").Value = 1673.125
$WVR.Range("J132").Value = 1942.3334
$WVR.Range("K132").Value = 5019.375
$WVR.Range("L132").Value = 5827.0002
$WVR.Range("M132").Value = -2489.375
$WVR.Range("N132").Value = -10887.0002

# WVR!row 136
$WVR.Range("H136").Value = 661.9583
$WVR.Range("I136").Value = 463.05884
$WVR.Range("J136").Value = 1145
$WVR.Range("K136").Value = 1389.17652
$WVR.Range("L136").Value = 3435
$WVR.Range("M136").Value = 1160.82348
$WVR.Range("N136").Value = -8535
